$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point drift in existing row 65's timestamp
$ws.Cells.Item(65, 1).Value = 44378.76918391898

# Append new row 66 with freshly retrieved data
$ws.Cells.Item(66, 1).Value = 44379.76524352474
$ws.Cells.Item(66, 2).Value = 78792
$ws.Cells.Item(66, 3).Value = 66394
$ws.Cells.Item(66, 4).Value = 3630
$ws.Cells.Item(66, 5).Value = 2124
$ws.Cells.Item(66, 6).Value = 1512
$ws.Cells.Item(66, 7).Value = 20929
$ws.Cells.Item(66, 8).Value = 1592
$ws.Cells.Item(66, 9).Value = 875
$ws.Cells.Item(66, 10).Value = 205

# New row's date cell should carry the same formatting as the rest of column A
$ws.Cells.Item(66, 1).NumberFormat = $ws.Cells.Item(65, 1).NumberFormat
